$d = $word.ActiveDocument

# --- 1. Redistribute the "tblEmployees" grid-column widths (twips) ---
# col1: 1339 -> 1338 ; col4: 1366 -> 1365 ; col7: 1360 -> 1362
$t = $d.Tables.Item(1)
$t.Columns.Item(1).Width = 1338 / 20.0
$t.Columns.Item(4).Width = 1365 / 20.0
$t.Columns.Item(7).Width = 1362 / 20.0

# --- 2. Remove the trailing "Thanks!" paragraph ---
# (table edits above can leave the Paragraphs collection's index cache
# stale, so locate the text via Find on an explicit Range instead of
# indexing into $d.Paragraphs)
$docEnd = $d.Content.End
$searchRng = $d.Range(0, $docEnd)
$found = $searchRng.Find.Execute("Thanks!")
if ($found) {
    $delRng = $d.Range($searchRng.Start, $docEnd)
    $delRng.Delete()
}

# --- 3. Register the two new "ListLabel" character styles ---
$s11 = $d.Styles.Add("ListLabel11", 2)
$s11.NameLocal = "ListLabel 11"
$s11.QuickStyle = $true
$s11.Font.Name = "Courier New"
$s11.Font.NameBi = "Courier New"
$s11.Font.Size = 9
$s11.Font.SizeBi = 9

$s12 = $d.Styles.Add("ListLabel12", 2)
$s12.NameLocal = "ListLabel 12"
$s12.QuickStyle = $true
$s12.Font.Name = "Courier New"
$s12.Font.NameBi = "Courier New"
$s12.Font.Size = 9
$s12.Font.SizeBi = 9
$s12.LanguageID = "en-US"

Write-Output "edit applied"
